$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Status" column before column D (Jan_2026 onwards shift right by one)
$ws.Columns.Item(4).Insert()

# Insert a new row before row 23 to make room for the extra data row
$ws.Rows.Item(23).Insert()

# --- Header row ---
$ws.Cells.Item(1,1).Value = "ISIN"
$ws.Cells.Item(1,2).Value = "Stock Name"
$ws.Cells.Item(1,3).Value = "Mutual Fund"
$ws.Cells.Item(1,4).Value = "Status"
$ws.Cells.Item(1,5).Value = "Jan_2026"
$ws.Cells.Item(1,6).Value = "Dec_2025"
$ws.Cells.Item(1,7).Value = "Oct_2025"
$ws.Cells.Item(1,8).Value = "MoM"
$ws.Cells.Item(1,9).Value = "QoQ"

# --- Data rows ---
# Row 2: HDFC Bank Limited
$ws.Cells.Item(2,1).Value = "INE040A01034"
$ws.Cells.Item(2,2).Value = "HDFC Bank Limited"
$ws.Cells.Item(2,3).Value = "quant Equity Savings Fund"
$ws.Cells.Item(2,4).Value = "Adding Consistently"
$ws.Cells.Item(2,5).Value = 9.269297
$ws.Cells.Item(2,6).Value = 5.688921
$ws.Cells.Item(2,7).Value = 5.747613
$ws.Cells.Item(2,8).Value = 3.580376
$ws.Cells.Item(2,9).Value = 3.521684

# Row 3: ICICI Bank Limited
$ws.Cells.Item(3,1).Value = "INE090A01021"
$ws.Cells.Item(3,2).Value = "ICICI Bank Limited"
$ws.Cells.Item(3,3).Value = "quant Equity Savings Fund"
$ws.Cells.Item(3,4).Value = "Adding Consistently"
$ws.Cells.Item(3,5).Value = 8.668876
$ws.Cells.Item(3,6).Value = 4.458872
$ws.Cells.Item(3,7).Value = 7.410043
$ws.Cells.Item(3,8).Value = 4.210003999999999
$ws.Cells.Item(3,9).Value = 1.258832999999999

# Row 4: Adani Green Energy Limited
$ws.Cells.Item(4,1).Value = "INE364U01010"
$ws.Cells.Item(4,2).Value = "Adani Green Energy Limited"
$ws.Cells.Item(4,3).Value = "quant Equity Savings Fund"
$ws.Cells.Item(4,4).Value = "Reducing"
$ws.Cells.Item(4,5).Value = 7.780966
$ws.Cells.Item(4,6).Value = 9.24471
$ws.Cells.Item(4,7).Value = 4.591021
$ws.Cells.Item(4,8).Value = -1.463743999999999
$ws.Cells.Item(4,9).Value = 3.189945000000001

# Row 5: Bharti Airtel Limited
$ws.Cells.Item(5,1).Value = "INE397D01024"
$ws.Cells.Item(5,2).Value = "Bharti Airtel Limited"
$ws.Cells.Item(5,3).Value = "quant Equity Savings Fund"
$ws.Cells.Item(5,4).Value = "Adding Consistently"
$ws.Cells.Item(5,5).Value = 6.797444
$ws.Cells.Item(5,6).Value = 5.218497
$ws.Cells.Item(5,7).Value = 5.164704
$ws.Cells.Item(5,8).Value = 1.578946999999999
$ws.Cells.Item(5,9).Value = 1.632739999999999

# Row 6: Lupin Limited
$ws.Cells.Item(6,1).Value = "INE326A01037"
$ws.Cells.Item(6,2).Value = "Lupin Limited"
$ws.Cells.Item(6,3).Value = "quant Equity Savings Fund"
$ws.Cells.Item(6,4).Value = "Adding Consistently"
$ws.Cells.Item(6,5).Value = 5.002997
$ws.Cells.Item(6,6).Value = 4.890458
$ws.Cells.Item(6,7).Value = 4.617115
$ws.Cells.Item(6,8).Value = 0.1125389999999999
$ws.Cells.Item(6,9).Value = 0.3858819999999996

# Row 7: Bajaj Finance Limited
$ws.Cells.Item(7,1).Value = "INE296A01032"
$ws.Cells.Item(7,2).Value = "Bajaj Finance Limited"
$ws.Cells.Item(7,3).Value = "quant Equity Savings Fund"
$ws.Cells.Item(7,4).Value = "Reducing Consistently"
$ws.Cells.Item(7,5).Value = 4.973993
$ws.Cells.Item(7,6).Value = 5.265807
$ws.Cells.Item(7,7).Value = 5.644254
$ws.Cells.Item(7,8).Value = -0.2918139999999996
$ws.Cells.Item(7,9).Value = -0.670261

# Row 8: Reliance Industries Limited
$ws.Cells.Item(8,1).Value = "INE002A01018"
$ws.Cells.Item(8,2).Value = "Reliance Industries Limited"
$ws.Cells.Item(8,3).Value = "quant Equity Savings Fund"
$ws.Cells.Item(8,4).Value = "Reducing Consistently"
$ws.Cells.Item(8,5).Value = 4.644474
$ws.Cells.Item(8,6).Value = 5.214247
$ws.Cells.Item(8,7).Value = 5.005953
$ws.Cells.Item(8,8).Value = -0.5697730000000005
$ws.Cells.Item(8,9).Value = -0.3614790000000001

# Row 9: Fortis Healthcare Ltd
$ws.Cells.Item(9,1).Value = "INE061F01013"
$ws.Cells.Item(9,2).Value = "Fortis Healthcare Ltd"
$ws.Cells.Item(9,3).Value = "quant Equity Savings Fund"
$ws.Cells.Item(9,4).Value = "Reducing Consistently"
$ws.Cells.Item(9,5).Value = 4.392929
$ws.Cells.Item(9,6).Value = 4.549517
$ws.Cells.Item(9,7).Value = 5.340472
$ws.Cells.Item(9,8).Value = -0.1565880000000002
$ws.Cells.Item(9,9).Value = -0.9475430000000005

# Row 10: Rural Electrification Corporation Ltd
$ws.Cells.Item(10,1).Value = "INE020B01018"
$ws.Cells.Item(10,2).Value = "Rural Electrification Corporation Ltd"
$ws.Cells.Item(10,3).Value = "quant Equity Savings Fund"
$ws.Cells.Item(10,4).Value = "Adding"
$ws.Cells.Item(10,5).Value = 4.362757
$ws.Cells.Item(10,6).Value = 4.264898
$ws.Cells.Item(10,7).Value = 4.483894
$ws.Cells.Item(10,8).Value = 0.09785900000000058
$ws.Cells.Item(10,9).Value = -0.1211370000000001

# Row 11: Kalyan Jewellers India Limited
$ws.Cells.Item(11,1).Value = "INE303R01014"
$ws.Cells.Item(11,2).Value = "Kalyan Jewellers India Limited"
$ws.Cells.Item(11,3).Value = "quant Equity Savings Fund"
$ws.Cells.Item(11,4).Value = "Reducing"
$ws.Cells.Item(11,5).Value = 4.344148
$ws.Cells.Item(11,6).Value = 5.815869
$ws.Cells.Item(11,7).Value = 0
$ws.Cells.Item(11,8).Value = -1.471721000000001
$ws.Cells.Item(11,9).Value = 4.344148

# Row 12: Eternal Limited
$ws.Cells.Item(12,1).Value = "INE758T01015"
$ws.Cells.Item(12,2).Value = "Eternal Limited"
$ws.Cells.Item(12,3).Value = "quant Equity Savings Fund"
$ws.Cells.Item(12,4).Value = "Reducing Consistently"
$ws.Cells.Item(12,5).Value = 4.258939
$ws.Cells.Item(12,6).Value = 4.317692
$ws.Cells.Item(12,7).Value = 5.00477
$ws.Cells.Item(12,8).Value = -0.05875300000000028
$ws.Cells.Item(12,9).Value = -0.7458309999999999

# Row 13: NMDC Ltd
$ws.Cells.Item(13,1).Value = "INE584A01023"
$ws.Cells.Item(13,2).Value = "NMDC Ltd"
$ws.Cells.Item(13,3).Value = "quant Equity Savings Fund"
$ws.Cells.Item(13,4).Value = "Reducing Consistently"
$ws.Cells.Item(13,5).Value = 1.954854
$ws.Cells.Item(13,6).Value = 1.99717
$ws.Cells.Item(13,7).Value = 4.922647
$ws.Cells.Item(13,8).Value = -0.0423159999999998
$ws.Cells.Item(13,9).Value = -2.967793

# Row 14: Kotak Mahindra Bank Limited
$ws.Cells.Item(14,1).Value = "INE237A01028"
$ws.Cells.Item(14,2).Value = "Kotak Mahindra Bank Limited"
$ws.Cells.Item(14,3).Value = "quant Equity Savings Fund"
$ws.Cells.Item(14,4).Value = "Complete Exit"
$ws.Cells.Item(14,5).Value = 0
$ws.Cells.Item(14,6).Value = 6.264325
$ws.Cells.Item(14,7).Value = 0
$ws.Cells.Item(14,8).Value = -6.264325
$ws.Cells.Item(14,9).Value = 0

# Row 15: JSW Steel Limited
$ws.Cells.Item(15,1).Value = "INE019A01038"
$ws.Cells.Item(15,2).Value = "JSW Steel Limited"
$ws.Cells.Item(15,3).Value = "quant Equity Savings Fund"
$ws.Cells.Item(15,4).Value = "Complete Exit"
$ws.Cells.Item(15,5).Value = 0
$ws.Cells.Item(15,6).Value = 0
$ws.Cells.Item(15,7).Value = 3.095238
$ws.Cells.Item(15,8).Value = 0
$ws.Cells.Item(15,9).Value = -3.095238

# Row 16: LIC Housing Finance Ltd
$ws.Cells.Item(16,1).Value = "INE115A01026"
$ws.Cells.Item(16,2).Value = "LIC Housing Finance Ltd"
$ws.Cells.Item(16,3).Value = "quant Equity Savings Fund"
$ws.Cells.Item(16,4).Value = "Complete Exit"
$ws.Cells.Item(16,5).Value = 0
$ws.Cells.Item(16,6).Value = 0
$ws.Cells.Item(16,7).Value = 2.919461
$ws.Cells.Item(16,8).Value = 0
$ws.Cells.Item(16,9).Value = -2.919461

# Row 17: Premier Energies Limited
$ws.Cells.Item(17,1).Value = "INE0BS701011"
$ws.Cells.Item(17,2).Value = "Premier Energies Limited"
$ws.Cells.Item(17,3).Value = "quant Equity Savings Fund"
$ws.Cells.Item(17,4).Value = "Complete Exit"
$ws.Cells.Item(17,5).Value = 0
$ws.Cells.Item(17,6).Value = 0
$ws.Cells.Item(17,7).Value = 5.950184
$ws.Cells.Item(17,8).Value = 0
$ws.Cells.Item(17,9).Value = -5.950184

# Row 18: State Bank of India
$ws.Cells.Item(18,1).Value = "INE062A01020"
$ws.Cells.Item(18,2).Value = "State Bank of India"
$ws.Cells.Item(18,3).Value = "quant Equity Savings Fund"
$ws.Cells.Item(18,4).Value = "Complete Exit"
$ws.Cells.Item(18,5).Value = 0
$ws.Cells.Item(18,6).Value = 0
$ws.Cells.Item(18,7).Value = 4.026852
$ws.Cells.Item(18,8).Value = 0
$ws.Cells.Item(18,9).Value = -4.026852

# Row 19: Apollo Hospitals Enterprise Ltd
$ws.Cells.Item(19,1).Value = "INE437A01024"
$ws.Cells.Item(19,2).Value = "Apollo Hospitals Enterprise Ltd"
$ws.Cells.Item(19,3).Value = "quant Equity Savings Fund"
$ws.Cells.Item(19,4).Value = "Complete Exit"
$ws.Cells.Item(19,5).Value = 0
$ws.Cells.Item(19,6).Value = 1.50322
$ws.Cells.Item(19,7).Value = 0
$ws.Cells.Item(19,8).Value = -1.50322
$ws.Cells.Item(19,9).Value = 0

# Row 20: Tata Consultancy Services Limited
$ws.Cells.Item(20,1).Value = "INE467B01029"
$ws.Cells.Item(20,2).Value = "Tata Consultancy Services Limited"
$ws.Cells.Item(20,3).Value = "quant Equity Savings Fund"
$ws.Cells.Item(20,4).Value = "Complete Exit"
$ws.Cells.Item(20,5).Value = 0
$ws.Cells.Item(20,6).Value = 4.92361
$ws.Cells.Item(20,7).Value = 4.763216
$ws.Cells.Item(20,8).Value = -4.92361
$ws.Cells.Item(20,9).Value = -4.763216

# Row 21: Godrej Properties Limited
$ws.Cells.Item(21,1).Value = "INE484J01027"
$ws.Cells.Item(21,2).Value = "Godrej Properties Limited"
$ws.Cells.Item(21,3).Value = "quant Equity Savings Fund"
$ws.Cells.Item(21,4).Value = "Complete Exit"
$ws.Cells.Item(21,5).Value = 0
$ws.Cells.Item(21,6).Value = 0
$ws.Cells.Item(21,7).Value = 2.820804
$ws.Cells.Item(21,8).Value = 0
$ws.Cells.Item(21,9).Value = -2.820804

# Row 22: SRF Limited
$ws.Cells.Item(22,1).Value = "INE647A01010"
$ws.Cells.Item(22,2).Value = "SRF Limited"
$ws.Cells.Item(22,3).Value = "quant Equity Savings Fund"
$ws.Cells.Item(22,4).Value = "Complete Exit"
$ws.Cells.Item(22,5).Value = 0
$ws.Cells.Item(22,6).Value = 0
$ws.Cells.Item(22,7).Value = 1.338015
$ws.Cells.Item(22,8).Value = 0
$ws.Cells.Item(22,9).Value = -1.338015

# Row 23: Bajaj Auto Limited
$ws.Cells.Item(23,1).Value = "INE917I01010"
$ws.Cells.Item(23,2).Value = "Bajaj Auto Limited"
$ws.Cells.Item(23,3).Value = "quant Equity Savings Fund"
$ws.Cells.Item(23,4).Value = "Complete Exit"
$ws.Cells.Item(23,5).Value = 0
$ws.Cells.Item(23,6).Value = 3.057868
$ws.Cells.Item(23,7).Value = 0
$ws.Cells.Item(23,8).Value = -3.057868
$ws.Cells.Item(23,9).Value = 0

